$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SSY")

$ws.Range("D9").Value = 18500
$ws.Range("D10").Value = 34300
$ws.Range("D17").Value = 54200
$ws.Range("D18").Value = -1300
$ws.Range("D20").Value = 200

$ws.Range("D21").Value = 700
$ws.Range("E21").Value = 200
$ws.Range("F21").Value = -4100
$ws.Range("G21").Value = 4900
$ws.Range("H21").Value = 4700
$ws.Range("I21").Value = 2200
$ws.Range("J21").Value = 6500

$ws.Range("D26").Value = -1100
$ws.Range("D27").Value = -1100
$ws.Range("D32").Value = -200
$ws.Range("D33").Value = -1600
$ws.Range("D35").Value = -1600
$ws.Range("D81").Value = -1600
